# Generate Report for Handoff
# Updates the localization-status report to reflect that the
# eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md file is now "Ready for handoff"
# (its handback xliff version is stale vs the latest source), across the
# Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e89bcb11f78cb4912ce5ef39800eaabe9585e374/e2e/eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a1d4611f7f423b81bcc4465fb008c6536abae77f/e2e/eaf3e711-3d40-44c3-a4d4-9772a4a9983d.md."

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-21 06:54:01"

# Excel stores column width internally as (ColumnWidth + 5/6); back the
# COM-visible ColumnWidth off by that padding so the saved <col width="..">
# lands on exactly 40.
$targetColWidth = 40 - (5 / 6)

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-21 06:53:56"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-21 06:54:01"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
